$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (post-edit) values for rows 5, 7, 8, 9, 10, 11, 12, 13, 15, 16, 17.
# Columns A, B, D, E, F, G, H, Q, R are the only ones that vary between rows;
# the rest of each row's data (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE,
# AG, AT, AW, AX, AY) is identical across all these rows and is left as-is.

$rows = @(
    @{ Row = 5;  A = 111943983; B = 90678;  D = "LC"; E = 4366;   F = "Skarp dropptaggsvamp";  G = "Hydnellum peckii";          H = "Banker";          Q = 682871.1304590552; R = 6694480.539619928 }
    @{ Row = 7;  A = 111943980; B = 89183;  D = "LC"; E = 3215;   F = "Rödgul trumpetsvamp";   G = "Craterellus lutescens";     H = "(Fr.) Fr.";       Q = 682877.1417635784; R = 6694410.432217407 }
    @{ Row = 8;  A = 111943992; B = 89183;  D = "LC"; E = 3215;   F = "Rödgul trumpetsvamp";   G = "Craterellus lutescens";     H = "(Fr.) Fr.";       Q = 682866.8554180798; R = 6694644.443727687 }
    @{ Row = 9;  A = 111943998; B = 98535;  D = "LC"; E = 222498; F = "Blåsippa";              G = "Hepatica nobilis";          H = "Schreb.";         Q = 682757.1772001419; R = 6694405.884787144 }
    @{ Row = 10; A = 111943990; B = 101703; D = "LC"; E = 222412; F = "Tibast";                G = "Daphne mezereum";           H = "L.";              Q = 682930.0967543643; R = 6694720.015570021 }
    @{ Row = 11; A = 111943995; B = 88899;  D = "NT"; E = 3286;   F = "Flattoppad klubbsvamp"; G = "Clavariadelphus truncatus"; H = "(Quél.) Donk";    Q = 682779.1674098044; R = 6694551.279700429 }
    @{ Row = 12; A = 111943988; B = 107033; D = "NT"; E = 220320; F = "Ängsskära";             G = "Serratula tinctoria";       H = "L.";              Q = 682930.0967543643; R = 6694720.015570021 }
    @{ Row = 13; A = 111943996; B = 90332;  D = "LC"; E = 4769;   F = "Svavelriska";           G = "Lactarius scrobiculatus";   H = "(Scop.:Fr.) Fr."; Q = 682785.3360249697; R = 6694547.127516991 }
    @{ Row = 15; A = 111943981; B = 96253;  D = "LC"; E = 504;    F = "Guckusko";              G = "Cypripedium calceolus";     H = "L.";              Q = 682877.1417635784; R = 6694410.432217407 }
    @{ Row = 16; A = 111943984; B = 99413;  D = "LC"; E = 221235; F = "Vårärt";                G = "Lathyrus vernus";           H = "(L.) Bernh.";     Q = 682929.3627028114; R = 6694685.271877083 }
    @{ Row = 17; A = 111943999; B = 99413;  D = "LC"; E = 221235; F = "Vårärt";                G = "Lathyrus vernus";           H = "(L.) Bernh.";     Q = 682757.1772001419; R = 6694405.884787144 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("Q$n").Value = $r.Q
    $ws.Range("R$n").Value = $r.R
}
